$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.726.37"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "1.860.14"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'304.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "'0.9984"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.5117"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("D8").Value = "'0.3649"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").Value = "'0.07135"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").Value = "'0.8879"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").Value = "'20.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "'0.07500"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "1.871.96"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "'94.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.56%  "
$ws.Range("D15").Value = "'5.209"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "'0.000008297"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.25%  "
$ws.Range("D18").Value = "'14.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "'0.9980"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "26.774.00"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "'4.985"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "2.111.69"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").Value = "'10.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "'6.351"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("D25").Value = "'145.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").Value = "'1.760"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.87%  "
$ws.Range("D27").Value = "'17.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "'2.074"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "'113.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "'4.672"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").Value = "'4.757"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").Value = "'0.09133"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.971"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.91%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7413"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("D36").Value = "'1.155"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "'3.217"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").Value = "'2.487"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01963"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.15%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5474"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("D41").Value = "'1.065"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").Value = "'6.504"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'115.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "'8.509"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("D45").Value = "'0.1480"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'0.4699"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").Value = "'0.9973"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").Value = "'9.990"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.547"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'36.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "'62.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.33%  "
